$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 606, shifting existing rows 606:686 down to 607:687
$ws.Rows.Item(606).Insert()

# Populate the new row 606 with the new daily price observation.
# (Columns A,B,C,E,F,G,H,I,J,Q,R,T repeat the same constant values used by
# every other row in this "Plátano / Vega Modelo de Temuco" table.)
$ws.Cells.Item(606, 1).Value = 10
$ws.Cells.Item(606, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(606, 3).Value = "La Araucanía"
$ws.Cells.Item(606, 4).Value = 44776
$ws.Cells.Item(606, 5).Value = 9
$ws.Cells.Item(606, 6).Value = "Fruta"
$ws.Cells.Item(606, 7).Value = 100108
$ws.Cells.Item(606, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(606, 9).Value = 100108006
$ws.Cells.Item(606, 10).Value = "Plátano"
$ws.Cells.Item(606, 11).Value = "Sin especificar"
$ws.Cells.Item(606, 12).Value = "Pintón"
$ws.Cells.Item(606, 13).Value = 150
$ws.Cells.Item(606, 14).Value = 33000
$ws.Cells.Item(606, 15).Value = 33000
$ws.Cells.Item(606, 16).Value = 33000
$ws.Cells.Item(606, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(606, 18).Value = "Ecuador"
$ws.Cells.Item(606, 19).Value = 1650
$ws.Cells.Item(606, 20).Value = 20
